$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value while forcing it to remain text (so numeric-looking
# strings like "247.19", "0.05600" or "-1.40%" keep their exact original formatting
# -- e.g. trailing zeros -- instead of Excel silently converting them to numbers).
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

Set-TextValue "D2" '247.19'
Set-TextValue "E2" '0.81%'
Set-TextValue "D3" '26.25'
Set-TextValue "E3" '4.52%'
Set-TextValue "D4" '5.084'
Set-TextValue "E4" '1.17%'
Set-TextValue "D5" '0.05600'
Set-TextValue "E5" '-0.36%'
Set-TextValue "D6" '6.475'
Set-TextValue "E6" '-1.46%'
Set-TextValue "D7" '0.8127'
Set-TextValue "E7" '-0.01%'
Set-TextValue "D8" '0.8446'
Set-TextValue "E8" '0.95%'
Set-TextValue "D9" '0.06984'
Set-TextValue "E9" '0.55%'
Set-TextValue "D10" '0.02809'
Set-TextValue "E10" '-1.14%'
Set-TextValue "D11" '0.09386'
Set-TextValue "E11" '-0.20%'
Set-TextValue "D12" '0.001512'
Set-TextValue "E12" '-1.19%'
Set-TextValue "D13" '0.0006009'
Set-TextValue "E13" '1.06%'
Set-TextValue "D14" '0.006148'
Set-TextValue "E14" '0.93%'
Set-TextValue "D15" '3.607'
Set-TextValue "E15" '3.09%'
Set-TextValue "D16" '3.019'
Set-TextValue "E16" '0.31%'
Set-TextValue "E17" '-1.74%'
Set-TextValue "D18" '0.3117'
Set-TextValue "E18" '-2.12%'
Set-TextValue "D19" '0.1339'
Set-TextValue "E19" '0.21%'
Set-TextValue "D20" '0.03203'
Set-TextValue "E20" '-1.85%'
Set-TextValue "D22" '3.747'
Set-TextValue "E22" '-0.09%'
Set-TextValue "D23" '0.04673'
Set-TextValue "E23" '0.20%'
Set-TextValue "E24" '-1.40%'
Set-TextValue "E25" '0.08%'
Set-TextValue "D27" '0.00009600'
Set-TextValue "E27" '-0.98%'
Set-TextValue "E28" '-0.04%'
Set-TextValue "D40" '0.03663'
Set-TextValue "E40" '0.01%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D41" '0.006130'
Set-TextValue "E41" '-1.55%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D42" '0.1056'
Set-TextValue "E42" '0.33%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D43" '0.002500'
Set-TextValue "E43" '-8.54%'
Set-TextValue "D44" '0.008291'
Set-TextValue "E44" '1.48%'
Set-TextValue "D45" '0.00005372'
Set-TextValue "E45" '1.48%'
Set-TextValue "E46" '0.06%'
Set-TextValue "E47" '-35.80%'
Set-TextValue "D48" '0.002583'
Set-TextValue "E48" '27.05%'
Set-TextValue "E49" '0.06%'
Set-TextValue "E50" '0.06%'
